$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add password123 for the access-level users (rows 3-6 / Alan, Scott, Paolo, Samuel)
$ws.Range("B3").Value = "password123"
$ws.Range("B4").Value = "password123"
$ws.Range("B5").Value = "password123"
$ws.Range("B6").Value = "password123"

# Adjust the access level for row 5 (Paolo Cisneros) from 3 to 2
$ws.Range("C5").Value = 2

# Widen column B slightly to fit the new password text (target stored width 12.5)
$ws.Range("B:B").ColumnWidth = 11.666666666666666

# Move the active selection to A5
$ws.Range("A5").Select()
